# Removed Test Case Inter-Dependency
$wb = $excel.ActiveWorkbook

$wsInput  = $wb.Worksheets.Item("ProductLoanInput")
$wsOutput = $wb.Worksheets.Item("ProductLoanOutput")

# productname: append "-1st" suffix on both sheets
$wsInput.Range("B1").Value  = "2490-RBI-EPP-DB-DL-REC-NOCOM-RNI-CTPD-SAR-MD-TR-2-DATE-VAR-INST-CASH-1st"
$wsOutput.Range("B1").Value = "2490-RBI-EPP-DB-DL-REC-NOCOM-RNI-CTPD-SAR-MD-TR-2-DATE-VAR-INST-CASH-1st"

# shortname: switch from numeric 2490 to text "249d"
$wsInput.Range("B2").Value = "249d"

# Selection on input sheet moves from B13 back to B1
$wsInput.Range("B1").Select()

# Output sheet becomes the active/selected tab
$wsOutput.Activate()
